# Qo'shimcha topshiriqlar / 07-amaliy ish: "yangi dasturlar va topshiriqlar"
#
# The "#N.[low,high] - N ball" lines lose their bold run/paragraph-mark
# formatting and get split into three runs ("#" / "N.[" / "low,high] - N ball")
# bracketed by <w:proofErr w:type="gramStart"/>...<w:proofErr w:type="gramEnd"/>
# around the middle run, exactly as Word's grammar checker would do once the
# bold toggle is removed. The two bare "# 5 ball" lines only lose their bold
# formatting (no text/run split).

$d = $word.ActiveDocument

# Common run/paragraph-mark properties shared by every affected paragraph,
# with the <w:b/><w:bCs/> pair dropped.
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="ru-RU"/></w:rPr>'
$pPr = '<w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/>' + $rPr + '</w:pPr>'
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Walk paragraphs back-to-front so inserting/replacing content never shifts
# the index of a paragraph we have not visited yet.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $bare = $text -replace "[\r\a]+$", ""

    if ($bare -match '^#(\d+)\.\[(.*)$') {
        # "#1.[0,10] - 1 ball" -> "#" | "1.[" | "0,10] - 1 ball"
        $num  = $matches[1]
        $rest = $matches[2]

        $xml = '<w:p' + $wNs + '>' + $pPr +
               '<w:r>' + $rPr + '<w:t>#</w:t></w:r>' +
               '<w:proofErr w:type="gramStart"/>' +
               '<w:r>' + $rPr + '<w:t>' + $num + '.[</w:t></w:r>' +
               '<w:proofErr w:type="gramEnd"/>' +
               '<w:r>' + $rPr + '<w:t>' + $rest + '</w:t></w:r>' +
               '</w:p>'

        [void]$para.Range.InsertXML($xml)
    }
    elseif ($bare -match '^#\s*5\s+ball$') {
        # "# 5 ball" keeps its single run/text, only bold goes away.
        $xml = '<w:p' + $wNs + '>' + $pPr +
               '<w:r>' + $rPr + '<w:t>' + $bare + '</w:t></w:r>' +
               '</w:p>'

        [void]$para.Range.InsertXML($xml)
    }
}
